# Add a new column "Serviced by " to the "Card22" worksheet, right after the
# existing "Correction" column, and tidy up the "Correction" header text
# (drop its trailing space). Also backfill the new/updated N column with the
# "nan" placeholder text used throughout the rest of the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card22")

# --- Header row (row 1) -----------------------------------------------
# N1: "Correction " -> "Correction" (trailing space removed)
$ws.Range("N1").Value = "Correction"

# O1: new header "Serviced by ", formatted the same as the other header
# cells (bold, centered, bordered) by copying N1's formatting.
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("O1").Value = "Serviced by "

# --- Data rows (rows 2-12) ---------------------------------------------
$lastRow = 12
for ($r = 2; $r -le $lastRow; $r++) {
    # N column: fill with "nan" like the rest of the row's columns
    $ws.Cells.Item($r, 14).Value = "nan"
    # O column: create the new (currently blank) cell for this row
    $ws.Cells.Item($r, 15).Font.Bold = $false
}
